$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 16, which was previously blank, with new log entry data.
$ws.Range("A16").Value = "Dovi IBS"
$ws.Range("B16").Value = "Lemon Shark"
$ws.Range("C16").Value = "Six"
$ws.Range("D5").Copy()
$ws.Range("D16").PasteSpecial(-4122)  # xlPasteFormats (adopt D5's date number format)
$ws.Range("D16").Value = "6/7/2021"
$ws.Range("E16").Value = "CKMR_DoviIBS_Lemon_sharks_AvgN_6yrs_06.07.2021_Lemon_CB_rmvTry.R"
$ws.Range("F16").Value = "Dovi_IBS_model_validation\Lemon_sharks\results\testing"
$ws.Range("G16").Value = 'Charlotte fixed the script! But she included a "try" section because errors were popping up occasionally. Here, I removed the "try" clause and saved each seed to see if the code runs and if it doesn''t, then I can use the seed to trace the error back and figure out what happened.'

# The wrapped text in G16 makes Excel auto-fit row 16's height (matches
# the other wrapped-text rows, e.g. row 8, which is also 75pt tall).
$ws.Rows.Item(16).RowHeight = 75

# Update the view: scroll so column D is the left-most visible column,
# and leave the active selection on H16.
$ws.Range("H16").Select()
$excel.ActiveWindow.ScrollColumn = 4
